$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3900
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(67, 8).Value = 3900
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(112, 8).Value = 1755.6428
$ws.Cells.Item(112, 9).Value = 790
$ws.Cells.Item(112, 10).Value = 1916.5834
$ws.Cells.Item(112, 11).Value = 2370
$ws.Cells.Item(112, 12).Value = 5749.7502
$ws.Cells.Item(112, 13).Value = -1262
$ws.Cells.Item(112, 14).Value = -7965.7502
$ws.Cells.Item(134, 8).Value = 32307.691
$ws.Cells.Item(134, 10).Value = 32307.691
$ws.Cells.Item(134, 12).Value = 32307.691
$ws.Cells.Item(134, 14).Value = -42447.691
$ws.Cells.Item(138, 8).Value = 1954.4259
$ws.Cells.Item(138, 9).Value = 1636.95
$ws.Cells.Item(138, 10).Value = 2141.1765
$ws.Cells.Item(138, 11).Value = 4910.85
$ws.Cells.Item(138, 12).Value = 6423.529500000001
$ws.Cells.Item(138, 13).Value = 229.1499999999996
$ws.Cells.Item(138, 14).Value = -16703.5295
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1488.1428
$ws.Cells.Item(2, 9).Value = 1159.3334
$ws.Cells.Item(2, 11).Value = 1159.3334
$ws.Cells.Item(2, 13).Value = -1046.3334
$ws.Cells.Item(45, 8).Value = 1321.0714
$ws.Cells.Item(45, 9).Value = 1255.5555
$ws.Cells.Item(45, 10).Value = 1439
$ws.Cells.Item(45, 11).Value = 1255.5555
$ws.Cells.Item(45, 12).Value = 1439
$ws.Cells.Item(45, 13).Value = -878.5554999999999
$ws.Cells.Item(45, 14).Value = -2193
$ws.Cells.Item(63, 8).Value = 590369.3
$ws.Cells.Item(63, 9).Value = 910987.2
$ws.Cells.Item(63, 10).Value = 2569.8333
$ws.Cells.Item(63, 11).Value = 910987.2
$ws.Cells.Item(63, 12).Value = 2569.8333
$ws.Cells.Item(63, 13).Value = -910301.2
$ws.Cells.Item(63, 14).Value = -3941.8333
$ws.Cells.Item(66, 8).Value = 590369.3
$ws.Cells.Item(66, 9).Value = 910987.2
$ws.Cells.Item(66, 10).Value = 2569.8333
$ws.Cells.Item(66, 11).Value = 4554936
$ws.Cells.Item(66, 12).Value = 12849.1665
$ws.Cells.Item(66, 13).Value = -4551504
$ws.Cells.Item(66, 14).Value = -19713.1665
$ws.Cells.Item(74, 8).Value = 1205.0741
$ws.Cells.Item(74, 9).Value = 1266.2174
$ws.Cells.Item(74, 10).Value = 853.5
$ws.Cells.Item(74, 11).Value = 1266.2174
$ws.Cells.Item(74, 12).Value = 853.5
$ws.Cells.Item(74, 13).Value = -392.2174
$ws.Cells.Item(74, 14).Value = -2601.5
$ws.Cells.Item(77, 8).Value = 1205.0741
$ws.Cells.Item(77, 9).Value = 1266.2174
$ws.Cells.Item(77, 10).Value = 853.5
$ws.Cells.Item(77, 11).Value = 6331.087
$ws.Cells.Item(77, 12).Value = 4267.5
$ws.Cells.Item(77, 13).Value = -1963.087
$ws.Cells.Item(77, 14).Value = -13003.5
$ws.Cells.Item(88, 8).Value = 2011403.8
$ws.Cells.Item(88, 9).Value = 2508003
$ws.Cells.Item(88, 10).Value = 25007
$ws.Cells.Item(88, 11).Value = 2508003
$ws.Cells.Item(88, 12).Value = 25007
$ws.Cells.Item(88, 13).Value = -2507597
$ws.Cells.Item(88, 14).Value = -25819
$ws.Cells.Item(91, 8).Value = 2011403.8
$ws.Cells.Item(91, 9).Value = 2508003
$ws.Cells.Item(91, 10).Value = 25007
$ws.Cells.Item(91, 11).Value = 2508003
$ws.Cells.Item(91, 12).Value = 25007
$ws.Cells.Item(91, 13).Value = -2506599
$ws.Cells.Item(91, 14).Value = -27815
$ws.Cells.Item(116, 8).Value = 1488.1428
$ws.Cells.Item(116, 9).Value = 1159.3334
$ws.Cells.Item(116, 11).Value = 1159.3334
$ws.Cells.Item(116, 13).Value = 1134.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1488.1428
$ws.Cells.Item(3, 9).Value = 1159.3334
$ws.Cells.Item(3, 11).Value = 1159.3334
$ws.Cells.Item(3, 13).Value = -1045.3334
$ws.Cells.Item(105, 8).Value = 5202
$ws.Cells.Item(105, 9).Value = 5202
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 5202
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -3455
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 28377.764
$ws.Cells.Item(134, 9).Value = 37466.605
$ws.Cells.Item(134, 10).Value = 2929
$ws.Cells.Item(134, 11).Value = 112399.815
$ws.Cells.Item(134, 12).Value = 8787
$ws.Cells.Item(134, 13).Value = -109864.815
$ws.Cells.Item(134, 14).Value = -13857
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3547828.2
$ws.Cells.Item(31, 9).Value = 1382.8158
$ws.Cells.Item(31, 10).Value = 18521710
$ws.Cells.Item(31, 11).Value = 1382.8158
$ws.Cells.Item(31, 12).Value = 18521710
$ws.Cells.Item(31, 13).Value = -1087.8158
$ws.Cells.Item(31, 14).Value = -18522300
$ws.Cells.Item(34, 8).Value = 3547828.2
$ws.Cells.Item(34, 9).Value = 1382.8158
$ws.Cells.Item(34, 10).Value = 18521710
$ws.Cells.Item(34, 11).Value = 1382.8158
$ws.Cells.Item(34, 12).Value = 18521710
$ws.Cells.Item(34, 13).Value = -1180.8158
$ws.Cells.Item(34, 14).Value = -18522114
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).ClearContents()
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).ClearContents()
$ws.Cells.Item(112, 8).Value = 19113.182
$ws.Cells.Item(112, 10).Value = 19113.182
$ws.Cells.Item(112, 12).Value = 19113.182
$ws.Cells.Item(112, 14).Value = -22067.182
$ws.Cells.Item(129, 8).Value = 49199.4
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 49199.4
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 49199.4
$ws.Cells.Item(129, 13).ClearContents()
$ws.Cells.Item(129, 14).Value = -59199.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 232.17949
$ws.Cells.Item(2, 9).Value = 309.7037
$ws.Cells.Item(2, 10).Value = 57.75
$ws.Cells.Item(2, 11).Value = 1858.2222
$ws.Cells.Item(2, 12).Value = 346.5
$ws.Cells.Item(2, 13).Value = -1745.2222
$ws.Cells.Item(2, 14).Value = -572.5
$ws.Cells.Item(132, 8).Value = 1003.73334
$ws.Cells.Item(132, 9).Value = 691.9583
$ws.Cells.Item(132, 11).Value = 6227.6247
$ws.Cells.Item(132, 13).Value = -3697.6247
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 7001
$ws.Cells.Item(80, 9).Value = 5000
$ws.Cells.Item(80, 10).Value = 7401.2
$ws.Cells.Item(80, 11).Value = 5000
$ws.Cells.Item(80, 12).Value = 7401.2
$ws.Cells.Item(80, 13).Value = -4002
$ws.Cells.Item(80, 14).Value = -9397.200000000001
$ws.Cells.Item(83, 8).Value = 7001
$ws.Cells.Item(83, 9).Value = 5000
$ws.Cells.Item(83, 10).Value = 7401.2
$ws.Cells.Item(83, 11).Value = 25000
$ws.Cells.Item(83, 12).Value = 37006
$ws.Cells.Item(83, 13).Value = -20008
$ws.Cells.Item(83, 14).Value = -46990
$ws.Cells.Item(129, 8).Value = 43332.832
$ws.Cells.Item(129, 10).Value = 43332.832
$ws.Cells.Item(129, 12).Value = 43332.832
$ws.Cells.Item(129, 14).Value = -53332.832
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1313.6666
$ws.Cells.Item(7, 10).Value = 2167.6667
$ws.Cells.Item(7, 12).Value = 2167.6667
$ws.Cells.Item(7, 14).Value = -2391.6667
$ws.Cells.Item(40, 8).Value = 2179
$ws.Cells.Item(40, 9).Value = 2179
$ws.Cells.Item(40, 11).Value = 2179
$ws.Cells.Item(40, 13).Value = -2043
$ws.Cells.Item(126, 8).Value = 1313.6666
$ws.Cells.Item(126, 10).Value = 2167.6667
$ws.Cells.Item(126, 12).Value = 6503.000100000001
$ws.Cells.Item(126, 14).Value = -11443.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 620.4375
$ws.Cells.Item(113, 9).Value = 654.9
$ws.Cells.Item(113, 10).Value = 563
$ws.Cells.Item(113, 11).Value = 1964.7
$ws.Cells.Item(113, 12).Value = 1689
$ws.Cells.Item(113, 13).Value = 205.3000000000002
$ws.Cells.Item(113, 14).Value = -6029
$ws.Cells.Item(126, 8).Value = 15025.375
$ws.Cells.Item(126, 9).Value = 27801
$ws.Cells.Item(126, 10).Value = 2249.75
$ws.Cells.Item(126, 11).Value = 83403
$ws.Cells.Item(126, 12).Value = 6749.25
$ws.Cells.Item(126, 13).Value = -80933
$ws.Cells.Item(126, 14).Value = -11689.25